# Applies the weekly refresh for 2023/england_league-one_2023-2024.xlsx
# - re-syncs the F:V (match) columns for rows whose underlying fixture order
#   shifted in the source feed (A:E identity columns are untouched)
# - appends the newly scraped fixture as row 186

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowFV {
    param($row, $vals)
    $arr = New-Object "object[,]" 1,17
    for ($i = 0; $i -lt 17; $i++) { $arr[0,$i] = $vals[$i] }
    $ws.Range("F$row`:V$row").Value = $arr
}

Set-RowFV 97 @("Burton", 1, "Fleetwood", 1, 2.2, "19/09/2023 20:12", 2.81, "23/09/2023 15:37", 3.41, "19/09/2023 20:12", 3.44, "23/09/2023 15:46", 3.38, "19/09/2023 20:12", 2.58, "23/09/2023 15:46", "https://www.betexplorer.com/football/england/league-one/burton-fleetwood-town/APxo5xLi/")
Set-RowFV 99 @("Carlisle", 0, "Derby", 2, 3.12, "19/09/2023 16:12", 3.53, "23/09/2023 15:35", 3.4, "19/09/2023 16:12", 3.43, "23/09/2023 15:35", 2.33, "19/09/2023 16:12", 2.18, "23/09/2023 15:35", "https://www.betexplorer.com/football/england/league-one/carlisle-derby/Iktg3Gk4/")
Set-RowFV 100 @("Bristol Rovers", 4, "Wigan", 1, 2.21, "19/09/2023 16:12", 2.23, "23/09/2023 15:48", 3.46, "19/09/2023 16:12", 3.6, "23/09/2023 15:37", 3.14, "19/09/2023 16:12", 3.26, "23/09/2023 15:48", "https://www.betexplorer.com/football/england/league-one/bristol-rovers-wigan/WYws6I5o/")
Set-RowFV 118 @("Charlton", 4, "Exeter", 1, 2.17, "30/09/2023 17:13", 1.87, "03/10/2023 20:35", 3.49, "30/09/2023 17:13", 3.83, "03/10/2023 20:35", 3.38, "30/09/2023 17:13", 4.17, "03/10/2023 20:35", "https://www.betexplorer.com/football/england/league-one/charlton-exeter/Ea21K8Qe/")
Set-RowFV 119 @("Carlisle", 1, "Peterborough", 1, 3.49, "30/09/2023 22:42", 3.66, "03/10/2023 20:44", 3.52, "30/09/2023 22:42", 3.62, "03/10/2023 20:44", 2.11, "30/09/2023 22:42", 2.06, "03/10/2023 20:44", "https://www.betexplorer.com/football/england/league-one/carlisle-peterborough/WzDdLlBk/")
Set-RowFV 120 @("Cheltenham", 0, "Fleetwood", 2, 2.98, "30/09/2023 17:13", 3.37, "03/10/2023 20:30", 3.32, "30/09/2023 17:13", 3.19, "03/10/2023 20:08", 2.36, "30/09/2023 17:13", 2.36, "03/10/2023 20:30", "https://www.betexplorer.com/football/england/league-one/cheltenham-fleetwood-town/hj15JSu2/")
Set-RowFV 121 @("Burton", 2, "Wigan", 1, 2.98, "30/09/2023 22:42", 2.62, "03/10/2023 20:31", 3.44, "30/09/2023 22:42", 3.28, "03/10/2023 20:22", 2.4, "30/09/2023 22:42", 2.88, "03/10/2023 20:31", "https://www.betexplorer.com/football/england/league-one/burton-wigan/2LMQQAmS/")
Set-RowFV 125 @("Cambridge Utd", 0, "Barnsley", 4, 2.61, "30/09/2023 17:13", 2.84, "03/10/2023 20:19", 3.37, "30/09/2023 17:13", 3.52, "03/10/2023 20:19", 2.64, "30/09/2023 17:13", 2.52, "03/10/2023 20:19", "https://www.betexplorer.com/football/england/league-one/cambridge-utd-barnsley/K09hMUeq/")
Set-RowFV 129 @("Peterborough", 2, "Lincoln", 0, 1.76, "03/10/2023 20:12", 1.61, "07/10/2023 15:52", 3.75, "03/10/2023 20:12", 4.07, "07/10/2023 16:00", 4.71, "03/10/2023 20:12", 5.88, "07/10/2023 16:00", "https://www.betexplorer.com/football/england/league-one/peterborough-lincoln-city/YyDfxsns/")
Set-RowFV 130 @("Leyton Orient", 2, "Reading", 1, 2.14, "03/10/2023 20:12", 2.4, "07/10/2023 15:57", 3.35, "03/10/2023 20:12", 3.35, "07/10/2023 15:52", 3.39, "03/10/2023 20:12", 3.14, "07/10/2023 15:52", "https://www.betexplorer.com/football/england/league-one/leyton-orient-reading/SxoNfk28/")
Set-RowFV 131 @("Portsmouth", 2, "Port Vale", 0, 1.66, "03/10/2023 20:12", 1.6, "07/10/2023 15:57", 3.95, "03/10/2023 20:12", 4.37, "07/10/2023 15:57", 5.22, "03/10/2023 20:12", 5.52, "07/10/2023 15:57", "https://www.betexplorer.com/football/england/league-one/portsmouth-port-vale/CSBby11m/")
Set-RowFV 132 @("Exeter", 0, "Barnsley", 1, 2.26, "03/10/2023 20:12", 2.93, "07/10/2023 15:52", 3.38, "03/10/2023 20:12", 3.46, "07/10/2023 15:52", 3.11, "03/10/2023 20:12", 2.48, "07/10/2023 15:52", "https://www.betexplorer.com/football/england/league-one/exeter-barnsley/CU2zXoA1/")
Set-RowFV 133 @("Cheltenham", 1, "Derby", 1, 4.29, "03/10/2023 20:12", 6.41, "07/10/2023 15:39", 3.7, "03/10/2023 20:12", 4.02, "07/10/2023 15:39", 1.78, "03/10/2023 20:12", 1.58, "07/10/2023 15:39", "https://www.betexplorer.com/football/england/league-one/cheltenham-derby/I73WXRfe/")
Set-RowFV 138 @("Blackpool", 3, "Stevenage", 0, 2.17, "07/10/2023 17:12", 2.18, "14/10/2023 15:58", 3.33, "07/10/2023 17:12", 3.35, "14/10/2023 15:29", 3.54, "07/10/2023 17:12", 3.57, "14/10/2023 15:58", "https://www.betexplorer.com/football/england/league-one/blackpool-stevenage/vu3FXI0C/")
Set-RowFV 139 @("Cambridge Utd", 1, "Shrewsbury", 1, 1.89, "09/10/2023 20:11", 2.04, "14/10/2023 15:58", 3.5, "09/10/2023 20:11", 3.29, "14/10/2023 15:58", 4.32, "09/10/2023 20:11", 4.15, "14/10/2023 15:58", "https://www.betexplorer.com/football/england/league-one/cambridge-utd-shrewsbury/z7dOVdVO/")
Set-RowFV 140 @("Carlisle", 0, "Leyton Orient", 1, 2.38, "07/10/2023 17:12", 2.72, "14/10/2023 15:54", 3.22, "07/10/2023 17:12", 3.34, "14/10/2023 15:54", 3.03, "07/10/2023 17:12", 2.73, "14/10/2023 15:52", "https://www.betexplorer.com/football/england/league-one/carlisle-leyton-orient/lAp6ft8P/")
Set-RowFV 141 @("Lincoln", 0, "Burton", 1, 1.85, "09/10/2023 20:11", 2.07, "14/10/2023 15:58", 3.6, "09/10/2023 20:11", 3.38, "14/10/2023 15:58", 4.09, "09/10/2023 20:11", 3.93, "14/10/2023 15:58", "https://www.betexplorer.com/football/england/league-one/lincoln-city-burton/8fRxkKVn/")
Set-RowFV 142 @("Leyton Orient", 1, "Barnsley", 1, 2.8, "14/10/2023 17:12", 2.67, "21/10/2023 15:46", 3.38, "14/10/2023 17:12", 3.35, "21/10/2023 14:59", 2.46, "14/10/2023 17:12", 2.78, "21/10/2023 15:46", "https://www.betexplorer.com/football/england/league-one/leyton-orient-barnsley/GvgSayp5/")
Set-RowFV 143 @("Stevenage", 0, "Port Vale", 0, 1.79, "14/10/2023 17:12", 1.97, "21/10/2023 15:51", 3.69, "14/10/2023 17:12", 3.5, "21/10/2023 15:52", 4.65, "14/10/2023 17:12", 4.15, "21/10/2023 15:51", "https://www.betexplorer.com/football/england/league-one/stevenage-port-vale/MuxBiDio/")
Set-RowFV 144 @("Shrewsbury", 1, "Derby", 0, 3.96, "14/10/2023 17:12", 4.39, "21/10/2023 15:56", 3.58, "14/10/2023 17:12", 3.52, "21/10/2023 15:56", 1.95, "14/10/2023 17:12", 1.91, "21/10/2023 15:42", "https://www.betexplorer.com/football/england/league-one/shrewsbury-derby/pKpqdgqU/")
Set-RowFV 145 @("Peterborough", 2, "Wycombe", 2, 1.75, "14/10/2023 05:13", 1.8, "21/10/2023 15:57", 3.88, "14/10/2023 05:13", 4.02, "21/10/2023 15:57", 4.63, "14/10/2023 05:13", 4.31, "21/10/2023 15:57", "https://www.betexplorer.com/football/england/league-one/peterborough-wycombe/O4dzbFEH/")
Set-RowFV 146 @("Oxford Utd", 1, "Blackpool", 1, 1.89, "14/10/2023 17:12", 2.03, "21/10/2023 15:56", 3.64, "14/10/2023 17:12", 3.53, "21/10/2023 15:56", 4.12, "14/10/2023 17:12", 3.85, "21/10/2023 15:56", "https://www.betexplorer.com/football/england/league-one/oxford-utd-blackpool/fViWbeaB/")
Set-RowFV 147 @("Fleetwood", 0, "Lincoln", 1, 2.19, "14/10/2023 17:12", 2.19, "21/10/2023 15:58", 3.32, "14/10/2023 17:12", 3.35, "21/10/2023 15:58", 3.51, "14/10/2023 17:12", 3.6, "21/10/2023 15:58", "https://www.betexplorer.com/football/england/league-one/fleetwood-town-lincoln-city/YF3N0HUb/")
Set-RowFV 148 @("Bolton", 2, "Northampton", 1, 1.65, "14/10/2023 05:13", 1.82, "21/10/2023 15:54", 4.06, "14/10/2023 05:13", 3.81, "21/10/2023 15:54", 4.72, "14/10/2023 05:13", 4.49, "21/10/2023 15:54", "https://www.betexplorer.com/football/england/league-one/bolton-northampton/OjH1qGaO/")
Set-RowFV 149 @("Charlton", 4, "Reading", 0, 1.83, "14/10/2023 13:43", 1.92, "21/10/2023 15:31", 3.77, "14/10/2023 13:43", 3.91, "21/10/2023 15:24", 3.96, "14/10/2023 13:43", 3.9, "21/10/2023 15:31", "https://www.betexplorer.com/football/england/league-one/charlton-reading/2HnC3Jpt/")
Set-RowFV 150 @("Burton", 4, "Bristol Rovers", 1, 2.76, "14/10/2023 17:12", 2.81, "21/10/2023 15:57", 3.34, "14/10/2023 17:12", 3.48, "21/10/2023 15:57", 2.62, "14/10/2023 17:12", 2.57, "21/10/2023 15:57", "https://www.betexplorer.com/football/england/league-one/burton-bristol-rovers/2sG5rzFU/")
Set-RowFV 151 @("Portsmouth", 1, "Carlisle", 0, 1.57, "14/10/2023 17:12", 1.42, "21/10/2023 15:19", 4.15, "14/10/2023 17:12", 4.85, "21/10/2023 15:19", 5.86, "14/10/2023 17:12", 7.77, "21/10/2023 15:19", "https://www.betexplorer.com/football/england/league-one/portsmouth-carlisle/2eevcZTN/")
Set-RowFV 152 @("Cheltenham", 1, "Cambridge Utd", 0, 2.77, "14/10/2023 17:12", 3.1, "21/10/2023 15:58", 3.14, "14/10/2023 17:12", 3.25, "21/10/2023 15:58", 2.75, "14/10/2023 17:12", 2.48, "21/10/2023 15:58", "https://www.betexplorer.com/football/england/league-one/cheltenham-cambridge-utd/t8oG2w0n/")
Set-RowFV 154 @("Wycombe", 2, "Bolton", 4, 2.88, "21/10/2023 17:13", 2.56, "24/10/2023 20:16", 3.36, "21/10/2023 17:13", 3.48, "24/10/2023 20:16", 2.51, "21/10/2023 17:13", 2.81, "24/10/2023 20:16", "https://www.betexplorer.com/football/england/league-one/wycombe-bolton/xYmRCrui/")
Set-RowFV 156 @("Port Vale", 0, "Peterborough", 1, 2.98, "21/10/2023 17:13", 3.31, "24/10/2023 20:29", 3.37, "21/10/2023 17:13", 3.6, "24/10/2023 20:29", 2.44, "21/10/2023 17:13", 2.2, "24/10/2023 20:29", "https://www.betexplorer.com/football/england/league-one/port-vale-peterborough/f9Ot5WLp/")
Set-RowFV 159 @("Lincoln", 3, "Charlton", 1, 2.65, "21/10/2023 17:13", 3.06, "24/10/2023 20:38", 3.39, "21/10/2023 17:13", 3.44, "24/10/2023 20:38", 2.58, "21/10/2023 17:13", 2.41, "24/10/2023 20:38", "https://www.betexplorer.com/football/england/league-one/lincoln-city-charlton/OQXAAhbN/")
Set-RowFV 166 @("Oxford Utd", 2, "Wycombe", 2, 1.71, "24/10/2023 21:12", 1.93, "28/10/2023 15:32", 3.88, "24/10/2023 21:12", 3.59, "28/10/2023 15:59", 4.91, "24/10/2023 21:12", 4.07, "28/10/2023 15:59", "https://www.betexplorer.com/football/england/league-one/oxford-utd-wycombe/p2aSmj87/")
Set-RowFV 168 @("Wigan", 2, "Shrewsbury", 0, 1.88, "24/10/2023 21:12", 2.17, "28/10/2023 15:57", 3.55, "24/10/2023 21:12", 3.32, "28/10/2023 15:57", 4.34, "24/10/2023 21:12", 3.69, "28/10/2023 15:57", "https://www.betexplorer.com/football/england/league-one/wigan-shrewsbury/YyeM9XUQ/")
Set-RowFV 169 @("Reading", 2, "Portsmouth", 3, 3.38, "24/10/2023 21:12", 3.56, "28/10/2023 15:56", 3.5, "24/10/2023 21:12", 3.57, "28/10/2023 15:56", 2.08, "24/10/2023 21:12", 2.11, "28/10/2023 15:56", "https://www.betexplorer.com/football/england/league-one/reading-portsmouth/xA1DBg0E/")
Set-RowFV 170 @("Port Vale", 1, "Cheltenham", 2, 1.75, "24/10/2023 21:12", 1.74, "28/10/2023 15:55", 3.63, "24/10/2023 21:12", 3.79, "28/10/2023 15:55", 4.58, "24/10/2023 21:12", 5.1, "28/10/2023 15:55", "https://www.betexplorer.com/football/england/league-one/port-vale-cheltenham/bJ29CZp8/")
Set-RowFV 172 @("Blackpool", 2, "Peterborough", 4, 2.13, "24/10/2023 21:12", 2.35, "28/10/2023 15:52", 3.46, "24/10/2023 21:12", 3.43, "28/10/2023 15:52", 3.31, "24/10/2023 21:12", 3.17, "28/10/2023 15:52", "https://www.betexplorer.com/football/england/league-one/blackpool-peterborough/MLkwAtP9/")
Set-RowFV 173 @("Cambridge Utd", 1, "Carlisle", 0, 2.04, "24/10/2023 21:12", 2.43, "28/10/2023 13:34", 3.13, "24/10/2023 21:12", 3.18, "28/10/2023 11:52", 3.96, "24/10/2023 21:12", 3.27, "28/10/2023 13:34", "https://www.betexplorer.com/football/england/league-one/cambridge-utd-carlisle/Aq7FjhOl/")
Set-RowFV 174 @("Burton", 0, "Leyton Orient", 0, 2.51, "24/10/2023 21:12", 3.16, "28/10/2023 15:55", 3.37, "24/10/2023 21:12", 3.37, "28/10/2023 15:56", 2.74, "24/10/2023 21:12", 2.38, "28/10/2023 15:56", "https://www.betexplorer.com/football/england/league-one/burton-leyton-orient/GM7BiY8r/")
Set-RowFV 175 @("Bristol Rovers", 2, "Northampton", 1, 1.95, "24/10/2023 21:12", 1.91, "28/10/2023 15:53", 3.59, "24/10/2023 21:12", 3.66, "28/10/2023 15:53", 3.95, "24/10/2023 21:12", 4.19, "28/10/2023 15:53", "https://www.betexplorer.com/football/england/league-one/bristol-rovers-northampton/vcpr90vG/")
Set-RowFV 183 @("Wigan", 2, "Peterborough", 1, 3.29, "05/11/2023 19:43", 3.07, "07/11/2023 20:43", 3.42, "05/11/2023 19:43", 3.76, "07/11/2023 20:43", 2.24, "05/11/2023 19:43", 2.27, "07/11/2023 20:42", "https://www.betexplorer.com/football/england/league-one/wigan-peterborough/GQDgoxVB/")
Set-RowFV 184 @("Shrewsbury", 0, "Bolton", 2, 3.58, "31/10/2023 07:11", 3.93, "07/11/2023 20:43", 3.37, "31/10/2023 07:11", 3.5, "07/11/2023 20:43", 2.14, "31/10/2023 07:11", 2.02, "07/11/2023 20:43", "https://www.betexplorer.com/football/england/league-one/shrewsbury-bolton/Qaa4d8xK/")

# Append the new fixture row (row 186), copying A/E number formats from row 185
$ws.Range("A185:V185").Copy()
$ws.Range("A186").PasteSpecial(-4122)

$ws.Range("A186").Value = 185
$ws.Range("B186").Value = "england"
$ws.Range("C186").Value = "league-one"
$ws.Range("D186").Value = "2023-2024"
$ws.Range("E186").Value = 45241.54166666666
Set-RowFV 186 @("Peterborough", 5, "Cambridge Utd", 0, 1.56, "03/11/2023 13:42", 1.51, "11/11/2023 12:53", 4.36, "03/11/2023 13:42", 4.44, "11/11/2023 12:53", 5.62, "03/11/2023 13:42", 6.72, "11/11/2023 12:53", "https://www.betexplorer.com/football/england/league-one/peterborough-cambridge-utd/t2KaNjaQ/")
